$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.985.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.108.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.372"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.870"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.03%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.106.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.689"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +20.96%  "

$ws.Range("E12").Value = "  +3.29%  "

$ws.Range("E13").Value = "  -5.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.034.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.709.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.108.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.08%  "

$ws.Range("E19").Value = "  +2.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000216"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "433.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "86.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.288.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.163"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "513.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.22%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.86%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.69%  "

$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.25%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.132"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.71%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +14.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.375"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.27%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0723"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +17.01%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "146.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "160.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.51%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
